$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (ACDA Tracking Sheet -> Heavy Planet (ACDA tracking v2))
$ws.Name = "Heavy Planet (ACDA tracking v2)"

# Keep the workbook's Print_Area / Print_Titles defined names in sync with
# the new sheet name (their scope prefix updates automatically, but the
# RefersTo text embeds the old sheet name and needs to be rewritten).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "='Heavy Planet (ACDA tracking v2)'!`$C`$3:`$O`$137"
    }
    if ($n.Name -like "*Print_Titles*") {
        $n.RefersTo = "='Heavy Planet (ACDA tracking v2)'!`$16:`$17"
    }
}

# Bump the sheet's zoom level from 100% to 110%
$excel.ActiveWindow.Zoom = 110

# Rows 47-51: L/M columns become literal 50 / 0 instead of the
# IF(L$24="","",L$24)-style lookup formula
foreach ($r in 47..51) {
    $ws.Range("L$r").Value = 50
    $ws.Range("M$r").Value = 0
}

# Row 54: N54 regains the IF(N$24="","",N$24) lookup formula it was missing
$ws.Range("N54").Formula = "=IF(N`$24=`"`",`"`",N`$24)"

# Rows 73, 78, 83, 88: L column changes from 0 to 50
foreach ($r in 73, 78, 83, 88) {
    $ws.Range("L$r").Value = 50
}

# Row 133: update several values
$ws.Range("F133").Value = 75
$ws.Range("G133").Value = 25
$ws.Range("H133").Value = 25
$ws.Range("I133").Value = 25
$ws.Range("J133").Value = 50
$ws.Range("L133").Value = 25
